# ---------------------------------------------------------------------------
# Renames "game lanes" terminology across the workbook and adds a localized
# "Restart Game" button label column to the Texts sheet, per the commit:
# "Changing spreadsheet path finding to release"
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename sheets
#    Trainings   -> Practices
#    Experiments -> Tasks
#    Prizes      -> Rewards
#    Times       -> Delays
# ---------------------------------------------------------------------------
$wsPractices = $wb.Worksheets.Item(1)
$wsPractices.Name = "Practices"

$wsTasks = $wb.Worksheets.Item(2)
$wsTasks.Name = "Tasks"

$wsTexts = $wb.Worksheets.Item(3)

$wsRewards = $wb.Worksheets.Item(4)
$wsRewards.Name = "Rewards"

$wsDelays = $wb.Worksheets.Item(5)
$wsDelays.Name = "Delays"

$wsDebug = $wb.Worksheets.Item(6)
$wsExport = $wb.Worksheets.Item(7)

# ---------------------------------------------------------------------------
# 2) Practices sheet header: "prize" -> "reward" terminology, second lane/value
#    renamed from "Second" to "Delayed"
# ---------------------------------------------------------------------------
$wsPractices.Range("B1").Value = "Immediate Reward Value"
$wsPractices.Range("C1").Value = "Delayed Reward Lane"
$wsPractices.Range("D1").Value = "Delayed Reward Value"

# ---------------------------------------------------------------------------
# 3) Tasks sheet: same header rename
# ---------------------------------------------------------------------------
$wsTasks.Range("B1").Value = "Immediate Reward Value"
$wsTasks.Range("C1").Value = "Delayed Reward Lane"
$wsTasks.Range("D1").Value = "Delayed Reward Value"

# ---------------------------------------------------------------------------
# 4) Texts sheet: insert a new "Restart Game" text column (X) between
#    "Task Score End" and "Score", with localized captions per language row.
# ---------------------------------------------------------------------------
$wsTexts.Columns.Item(24).Insert()
$wsTexts.Range("X1").Value = "Restart Game"
$wsTexts.Range("X2").Value = "Reiniciar o jogo"
$wsTexts.Range("X3").Value = "Restart game"
$wsTexts.Range("X4").Value = "Reiniciar el juego"

# ---------------------------------------------------------------------------
# 5) Delays sheet header: "pista"/"tempo" -> "Lane"/"Time"
# ---------------------------------------------------------------------------
$wsDelays.Range("A1").Value = "Lane"
$wsDelays.Range("B1").Value = "Time"

# ---------------------------------------------------------------------------
# 6) Restore per-sheet selections / active cells, and move the active tab
#    from Export back to Practices.
# ---------------------------------------------------------------------------
$wsTasks.Activate()
$wsTasks.Range("G6").Select()

$wsTexts.Activate()
$wsTexts.Range("X11").Select()

$wsDelays.Activate()
$wsDelays.Range("G9").Select()

$wsDebug.Activate()
$wsDebug.Range("C3").Select()

$wsExport.Activate()
$wsExport.Range("E14").Select()

$wsPractices.Activate()
$wsPractices.Range("E16").Select()
